$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shared-string table reorder shifted the labels that share rows 4-6:
# B4 now reads "Mkm_s2l", B5 now reads "Cpu_s2l", B6 now reads "EDN_S2l".
$ws.Range("B4").Value = "Mkm_s2l"
$ws.Range("B5").Value = "Cpu_s2l"
$ws.Range("B6").Value = "EDN_S2l"

# Updated Weekly Pending Total(Rp) and Repayment figures (column C / D).
$ws.Range("C2").Value = 7233217270
$ws.Range("D2").Value = 831371019

$ws.Range("C3").Value = 2985583999
$ws.Range("D3").Value = 337219452

$ws.Range("C4").Value = 3893666823
$ws.Range("D4").Value = 417201143

$ws.Range("C5").Value = 1991255199
$ws.Range("D5").Value = 212572570

$ws.Range("C6").Value = 1981408169
$ws.Range("D6").Value = 205479337

$ws.Range("C7").Value = 7356278731
$ws.Range("D7").Value = 735906523

$ws.Range("C8").Value = 7173980429
$ws.Range("D8").Value = 692573827
